$wb = $excel.ActiveWorkbook

# Template sheet that Denmark/Sweden/Norway are cloned from ("UK" tab)
$template = $wb.Worksheets.Item("UK")

# --- Denmark --------------------------------------------------------------
$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$denmark = $wb.Worksheets.Item($wb.Worksheets.Count)
$denmark.Name = "Denmark"
$denmark.Range("B2").Value = "Denmark Market"
$denmark.Range("B4").Value = "NGC-3446/T2004/T2005"

# --- Sweden -----------------------------------------------------------------
$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$sweden = $wb.Worksheets.Item($wb.Worksheets.Count)
$sweden.Name = "Sweden"
$sweden.Range("B2").Value = "Sweden Market"
$sweden.Range("B4").Value = "NGC-3465/T2025/T2023"

# --- Norway -----------------------------------------------------------------
$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$norway = $wb.Worksheets.Item($wb.Worksheets.Count)
$norway.Name = "Norway"
$norway.Range("B2").Value = "Norway Market"
$norway.Range("B4").Value = "NGC-3464/T1919"

# Norway ends up the active/selected tab, with B2:B4 selected (anchor B2)
$norway.Activate()
$norway.Range("B2:B4").Select()
